# Auto update Excel log
# Appends new sensor/alert rows to the ALERTS, Humidity, Temperature and
# Proximity sheets, matching the latest sensor sweep recorded on 2026-02-01.

function Set-LogRow {
    param($ws, $r, $date, $timestamp, $hour, $location, $value, $status)
    # Dates such as "2026-02-01" are quote-prefixed so Excel keeps them as
    # literal text instead of silently converting them to date serials.
    $ws.Range("A$r").Value = "'" + $date
    $ws.Range("B$r").Value = $timestamp
    $ws.Range("C$r").Value = $hour
    $ws.Range("D$r").Value = $location
    # Percentage-looking values (e.g. "78.8%") also need the quote prefix so
    # they remain plain text rather than becoming numeric percentages.
    if ($value -match '^[0-9.]+%$') {
        $ws.Range("E$r").Value = "'" + $value
    } else {
        $ws.Range("E$r").Value = $value
    }
    $ws.Range("F$r").Value = $status
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALERTS sheet: append row 8 (MINIMAL bathroom alert)
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Set-LogRow $wsAlerts 8 "2026-02-01" "18:25:00" "18:00" "Bathroom" "MINIMAL" "MINIMAL ALERT: Bathroom occupied, no motion > 20s."

# ---------------------------------------------------------------------
# Humidity sheet: append rows 65-73
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
Set-LogRow $wsHumidity 65 "2026-02-01" "18:24:18" "18:00" "Bathroom" "78.8%" "Active"
Set-LogRow $wsHumidity 66 "2026-02-01" "18:24:32" "18:00" "Bathroom" "78.9%" "Active"
Set-LogRow $wsHumidity 67 "2026-02-01" "18:24:37" "18:00" "Bathroom" "80.5%" "Active"
Set-LogRow $wsHumidity 68 "2026-02-01" "18:24:42" "18:00" "Bathroom" "88.1%" "Active"
Set-LogRow $wsHumidity 69 "2026-02-01" "18:24:47" "18:00" "Bathroom" "94.4%" "Active"
Set-LogRow $wsHumidity 70 "2026-02-01" "18:24:52" "18:00" "Bathroom" "87.5%" "Active"
Set-LogRow $wsHumidity 71 "2026-02-01" "18:24:57" "18:00" "Bathroom" "90.2%" "Active"
Set-LogRow $wsHumidity 72 "2026-02-01" "18:25:13" "18:00" "Bathroom" "83.1%" "Active"
Set-LogRow $wsHumidity 73 "2026-02-01" "18:25:18" "18:00" "Bathroom" "82.1%" "Active"

# ---------------------------------------------------------------------
# Temperature sheet: append rows 65-73
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
Set-LogRow $wsTemperature 65 "2026-02-01" "18:24:19" "18:00" "Bathroom" "29.4C" "Active"
Set-LogRow $wsTemperature 66 "2026-02-01" "18:24:33" "18:00" "Bathroom" "29.4C" "Active"
Set-LogRow $wsTemperature 67 "2026-02-01" "18:24:38" "18:00" "Bathroom" "29.4C" "Active"
Set-LogRow $wsTemperature 68 "2026-02-01" "18:24:43" "18:00" "Bathroom" "29.4C" "Active"
Set-LogRow $wsTemperature 69 "2026-02-01" "18:24:48" "18:00" "Bathroom" "29.5C" "Active"
Set-LogRow $wsTemperature 70 "2026-02-01" "18:24:53" "18:00" "Bathroom" "29.6C" "Active"
Set-LogRow $wsTemperature 71 "2026-02-01" "18:24:58" "18:00" "Bathroom" "29.6C" "Active"
Set-LogRow $wsTemperature 72 "2026-02-01" "18:25:13" "18:00" "Bathroom" "29.7C" "Active"
Set-LogRow $wsTemperature 73 "2026-02-01" "18:25:18" "18:00" "Bathroom" "29.7C" "Active"

# ---------------------------------------------------------------------
# Proximity sheet: append row 42 (bathroom door ENTER event)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
Set-LogRow $wsProximity 42 "2026-02-01" "18:24:37" "18:00" "Bathroom Door" "ENTER" "User ENTERED Bathroom"
